$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 3 title text (was a "TITLE" placeholder) and the newly-recorded times.
$ws.Range("B7").Value = "Day 3: Crossed Wires"
$ws.Range("C7").Value = 0.015381944444444443
$ws.Range("E7").Value = 0.020300925925925927
$ws.Range("F7").Value = 0.010902777777777777
$ws.Range("H7").Value = "2nd"

# Move the active selection down to H8, matching where entry continued.
$null = $ws.Range("H8").Select()
